$d = $word.ActiveDocument

# Locate the paragraph that ends the bibliography ("...GIL, A.C. ... 2010.").
# The three paragraphs that follow it -- a blank paragraph, the
# "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph, and the
# "(c) 2020 . Contact: ..." paragraph -- are removed, while the blank
# paragraph and page-break paragraph that come after them are kept.

$gilPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*2005. GIL, A.C. Como elaborar projetos de pesquisa*") {
        $gilPara = $p
        break
    }
}

$p1 = $gilPara.Next()
$p2 = $p1.Next()
$p3 = $p2.Next()

$deleteRange = $d.Range($p1.Range.Start, $p3.Range.End)
$deleteRange.Delete()
